$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 426.66666
$ws.Range("I6").Value = 426.66666
$ws.Range("K6").Value = 1279.99998
$ws.Range("M6").Value = -1167.99998
$ws.Range("H41").Value = 251.85715
$ws.Range("I41").Value = 200
$ws.Range("J41").Value = 321
$ws.Range("K41").Value = 200
$ws.Range("L41").Value = 321
$ws.Range("M41").Value = 240
$ws.Range("N41").Value = -1201
$ws.Range("H55").Value = 520.6667
$ws.Range("I55").Value = 350
$ws.Range("J55").Value = 569.4286
$ws.Range("K55").Value = 350
$ws.Range("L55").Value = 569.4286
$ws.Range("M55").Value = -136
$ws.Range("N55").Value = -997.4286
$ws.Range("H98").Value = 1830.1364
$ws.Range("I98").Value = 1498.3125
$ws.Range("J98").Value = 2715
$ws.Range("K98").Value = 1498.3125
$ws.Range("L98").Value = 2715
$ws.Range("M98").Value = -0.3125
$ws.Range("N98").Value = -5711
$ws.Range("H122").Value = 1830.1364
$ws.Range("I122").Value = 1498.3125
$ws.Range("J122").Value = 2715
$ws.Range("K122").Value = 4494.9375
$ws.Range("L122").Value = 8145
$ws.Range("M122").Value = -2044.9375
$ws.Range("N122").Value = -13045
$ws.Range("H125").Value = 5809.364
$ws.Range("I125").Value = 1110.6666
$ws.Range("J125").Value = 6551.263
$ws.Range("K125").Value = 9995.999400000001
$ws.Range("L125").Value = 58961.367
$ws.Range("M125").Value = -7535.999400000001
$ws.Range("N125").Value = -63881.367
$ws.Range("H132").Value = 1276.1305
$ws.Range("I132").Value = 1276.1305
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3828.3915
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1298.3915
$ws.Range("N132").ClearContents()
$ws.Range("H137").Value = 697293.4
$ws.Range("I137").Value = 3883.3333
$ws.Range("J137").Value = 1012479.75
$ws.Range("K137").Value = 11649.9999
$ws.Range("L137").Value = 3037439.25
$ws.Range("M137").Value = -9099.999899999999
$ws.Range("N137").Value = -3042539.25
$ws.Range("H138").Value = 4786.6113
$ws.Range("I138").Value = 2518.4546
$ws.Range("J138").Value = 5366.8374
$ws.Range("K138").Value = 7555.3638
$ws.Range("L138").Value = 16100.5122
$ws.Range("M138").Value = -2415.3638
$ws.Range("N138").Value = -26380.5122

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19153.678
$ws.Range("I32").Value = 19875.861
$ws.Range("J32").Value = 8682
$ws.Range("K32").Value = 19875.861
$ws.Range("L32").Value = 8682
$ws.Range("M32").Value = -19588.861
$ws.Range("N32").Value = -9256
$ws.Range("H41").Value = 7991.2
$ws.Range("I41").Value = 4978
$ws.Range("K41").Value = 4978
$ws.Range("M41").Value = -4564
$ws.Range("H45").Value = 2122.2727
$ws.Range("I45").Value = 2034.5
$ws.Range("K45").Value = 2034.5
$ws.Range("M45").Value = -1657.5
$ws.Range("H75").Value = 40173
$ws.Range("J75").Value = 40173
$ws.Range("L75").Value = 40173
$ws.Range("N75").Value = -41921
$ws.Range("H78").Value = 40173
$ws.Range("J78").Value = 40173
$ws.Range("L78").Value = 120519
$ws.Range("N78").Value = -129255
$ws.Range("H109").Value = 60063.5
$ws.Range("J109").Value = 60063.5
$ws.Range("L109").Value = 60063.5
$ws.Range("N109").Value = -62837.5
$ws.Range("H123").Value = 78900
$ws.Range("J123").Value = 78900
$ws.Range("L123").Value = 78900
$ws.Range("N123").Value = -88700

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 2999
$ws.Range("I24").Value = 2332
$ws.Range("K24").Value = 2332
$ws.Range("M24").Value = -2097
$ws.Range("H25").Value = 2271.3333
$ws.Range("I25").Value = 2271.3333
$ws.Range("K25").Value = 2271.3333
$ws.Range("M25").Value = -2036.3333

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 450096.47
$ws.Range("I31").Value = 5210.5815
$ws.Range("J31").Value = 928348.75
$ws.Range("K31").Value = 5210.5815
$ws.Range("L31").Value = 928348.75
$ws.Range("M31").Value = -4915.5815
$ws.Range("N31").Value = -928938.75
$ws.Range("H34").Value = 450096.47
$ws.Range("I34").Value = 5210.5815
$ws.Range("J34").Value = 928348.75
$ws.Range("K34").Value = 5210.5815
$ws.Range("L34").Value = 928348.75
$ws.Range("M34").Value = -5008.5815
$ws.Range("N34").Value = -928752.75
$ws.Range("H60").Value = 22500
$ws.Range("J60").Value = 22500
$ws.Range("L60").Value = 22500
$ws.Range("N60").Value = -23522
$ws.Range("H99").Value = 2091.6667
$ws.Range("I99").Value = 1614.2858
$ws.Range("J99").Value = 2760
$ws.Range("K99").Value = 1614.2858
$ws.Range("L99").Value = 2760
$ws.Range("M99").Value = -116.2858000000001
$ws.Range("N99").Value = -5756
$ws.Range("H126").Value = 2091.6667
$ws.Range("I126").Value = 1614.2858
$ws.Range("J126").Value = 2760
$ws.Range("K126").Value = 4842.857400000001
$ws.Range("L126").Value = 8280
$ws.Range("M126").Value = -2372.857400000001
$ws.Range("N126").Value = -13220
$ws.Range("H133").Value = 25886.375
$ws.Range("I133").Value = 10098.667
$ws.Range("J133").Value = 35359
$ws.Range("K133").Value = 10098.667
$ws.Range("L133").Value = 35359
$ws.Range("M133").Value = -7568.666999999999
$ws.Range("N133").Value = -40419
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 401
$ws.Range("I16").Value = 401
$ws.Range("K16").Value = 1203
$ws.Range("M16").Value = -1030
$ws.Range("H47").Value = 744.375
$ws.Range("I47").Value = 640.6
$ws.Range("J47").Value = 917.3333
$ws.Range("K47").Value = 1921.8
$ws.Range("L47").Value = 2751.9999
$ws.Range("M47").Value = -1490.8
$ws.Range("N47").Value = -3613.9999
$ws.Range("H63").Value = 4849.5
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 4849.5
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 14548.5
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -16046.5
$ws.Range("H66").Value = 4849.5
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 4849.5
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 43645.5
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -51133.5
$ws.Range("H68").Value = 1609.3066
$ws.Range("I68").Value = 1443.2
$ws.Range("J68").Value = 1858.4667
$ws.Range("K68").Value = 4329.6
$ws.Range("L68").Value = 5575.4001
$ws.Range("M68").Value = -3518.6
$ws.Range("N68").Value = -7197.4001
$ws.Range("H71").Value = 1609.3066
$ws.Range("I71").Value = 1443.2
$ws.Range("J71").Value = 1858.4667
$ws.Range("K71").Value = 12988.8
$ws.Range("L71").Value = 16726.2003
$ws.Range("M71").Value = -8932.800000000001
$ws.Range("N71").Value = -24838.2003

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3309.0334
$ws.Range("I102").Value = 2619.3333
$ws.Range("J102").Value = 4343.5835
$ws.Range("K102").Value = 2619.3333
$ws.Range("L102").Value = 4343.5835
$ws.Range("M102").Value = -997.3332999999998
$ws.Range("N102").Value = -7587.5835
$ws.Range("H135").Value = 61356
$ws.Range("J135").Value = 61356
$ws.Range("L135").Value = 61356
$ws.Range("N135").Value = -71496

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 10000
$ws.Range("I13").Value = 10000
$ws.Range("K13").Value = 10000
$ws.Range("M13").Value = -9860
$ws.Range("H16").Value = 935.4
$ws.Range("I16").Value = 944.25
$ws.Range("J16").Value = 900
$ws.Range("K16").Value = 944.25
$ws.Range("L16").Value = 900
$ws.Range("M16").Value = -774.25
$ws.Range("N16").Value = -1240
$ws.Range("H69").Value = 33358774
$ws.Range("J69").Value = 33358774
$ws.Range("L69").Value = 33358774
$ws.Range("N69").Value = -33360396
$ws.Range("H72").Value = 33358774
$ws.Range("J72").Value = 33358774
$ws.Range("L72").Value = 100076322
$ws.Range("N72").Value = -100084434
$ws.Range("H93").Value = 985
$ws.Range("I93").Value = 999.2308
$ws.Range("J93").Value = 800
$ws.Range("K93").Value = 999.2308
$ws.Range("L93").Value = 800
$ws.Range("M93").Value = 248.7692
$ws.Range("N93").Value = -3296
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 12852.381
$ws.Range("I15").Value = 9993.333000000001
$ws.Range("K15").Value = 9993.333000000001
$ws.Range("M15").Value = -9705.333000000001
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H110").Value = 41000
$ws.Range("I110").Value = 30000
$ws.Range("J110").Value = 52000
$ws.Range("K110").Value = 30000
$ws.Range("L110").Value = 52000
$ws.Range("M110").Value = -25910
$ws.Range("N110").Value = -60180
$ws.Range("H123").Value = 34862.1
$ws.Range("J123").Value = 34862.1
$ws.Range("L123").Value = 34862.1
$ws.Range("N123").Value = -44662.1
$ws.Range("H126").Value = 1821.2106
$ws.Range("I126").Value = 1787.6875
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 5363.0625
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -2893.0625
$ws.Range("N126").Value = -10940
